$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.32"
$ws.Range("E2").Value = "'-1.05%"
$ws.Range("D3").Value = "'43.95"
$ws.Range("E3").Value = "'5.83%"
$ws.Range("D4").Value = "'5.495"
$ws.Range("E4").Value = "'-3.38%"
$ws.Range("D5").Value = "'0.08068"
$ws.Range("E5").Value = "'-4.25%"
$ws.Range("D6").Value = "'8.633"
$ws.Range("E6").Value = "'-2.03%"
$ws.Range("D7").Value = "'4.276"
$ws.Range("E7").Value = "'-4.71%"
$ws.Range("D8").Value = "'1.881"
$ws.Range("E8").Value = "'-5.55%"
$ws.Range("D10").Value = "'0.9355"
$ws.Range("E10").Value = "'0.90%"
$ws.Range("D11").Value = "'0.1158"
$ws.Range("E11").Value = "'-9.12%"
$ws.Range("D12").Value = "'0.1892"
$ws.Range("E12").Value = "'-3.62%"
$ws.Range("D13").Value = "'0.09617"
$ws.Range("E13").Value = "'2.87%"
$ws.Range("D14").Value = "'0.04165"
$ws.Range("E14").Value = "'5.28%"
$ws.Range("D15").Value = "'0.1065"
$ws.Range("E15").Value = "'0.23%"
$ws.Range("D16").Value = "'0.001273"
$ws.Range("E16").Value = "'-2.24%"
$ws.Range("D17").Value = "'0.005981"
$ws.Range("E17").Value = "'-2.21%"
$ws.Range("D18").Value = "'3.574"
$ws.Range("E18").Value = "'4.09%"
$ws.Range("D20").Value = "'8.562"
$ws.Range("E20").Value = "'-6.75%"
$ws.Range("E21").Value = "'-0.01%"
$ws.Range("E22").Value = "'3.25%"
$ws.Range("D23").Value = "'0.04341"
$ws.Range("E23").Value = "'-1.72%"
$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'-1.07%"
$ws.Range("D25").Value = "'0.004392"
$ws.Range("E25").Value = "'-0.18%"
$ws.Range("D26").Value = "'0.0001233"
$ws.Range("E26").Value = "'3.47%"
$ws.Range("D27").Value = "'0.0003996"
$ws.Range("E27").Value = "'-0.01%"
$ws.Range("D39").Value = "'0.02648"
$ws.Range("E39").Value = "'-6.91%"
$ws.Range("D40").Value = "'0.05450"
$ws.Range("E40").Value = "'-1.44%"
$ws.Range("D41").Value = "'0.01144"
$ws.Range("E41").Value = "'27.44%"
$ws.Range("D42").Value = "'0.007649"
$ws.Range("E42").Value = "'-3.16%"
$ws.Range("D43").Value = "'0.1388"
$ws.Range("E43").Value = "'-3.37%"
$ws.Range("D44").Value = "'0.002112"
$ws.Range("E44").Value = "'1.41%"
$ws.Range("D45").Value = "'0.009626"
$ws.Range("E45").Value = "'-12.34%"
$ws.Range("D46").Value = "'0.00006913"
$ws.Range("E46").Value = "'-4.71%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("D48").Value = "'0.003565"
$ws.Range("E48").Value = "'10.57%"
$ws.Range("D49").Value = "'0.002274"
$ws.Range("E49").Value = "'-0.34%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.01%"

$ws.Range("D2:E51").Style = "Normal"
